{"js": "// Update the date heading and the 25 multiplication problems laid out in\n// the single 20-row x 5-column table (only rows 0, 4, 9, 14, 19 hold the\n// actual \"NNN\u00d7N=\" problems; the other rows are blank spacer rows).\n//\n// Cell addressing (row, col) is used instead of text search-and-replace so\n// that the several repeated/aliased values (e.g. one problem's new value\n// equals another problem's old value) can never be mismatched.\n\nconst body = context.document.body;\n\n// 1) Date heading, first paragraph of the body.\nconst headingPara = body.paragraphs.getFirst();\nheadingPara.load(\"text\");\nawait context.sync();\n\nif (headingPara.text.trim() === \"2025-03-02 Sunday\") {\n  headingPara.getRange().insertText(\"2025-03-03 Monday\", \"Replace\");\n  await context.sync();\n}\n\n// 2) Table of problems.\nconst table = body.tables.getFirst();\n\nconst rowValues = {\n  0: [\"916\u00d79=\", \"568\u00d74=\", \"489\u00d77=\", \"965\u00d78=\", \"525\u00d78=\"],\n  4: [\"366\u00d78=\", \"964\u00d72=\", \"544\u00d74=\", \"985\u00d73=\", \"724\u00d72=\"],\n  9: [\"594\u00d78=\", \"323\u00d76=\", \"728\u00d72=\", \"714\u00d74=\", \"343\u00d78=\"],\n  14: [\"809\u00d77=\", \"618\u00d74=\", \"805\u00d75=\", \"181\u00d74=\", \"983\u00d79=\"],\n  19: [\"800\u00d73=\", \"866\u00d79=\", \"746\u00d79=\", \"731\u00d78=\", \"201\u00d79=\"],\n};\n\nfor (const rowIndex of Object.keys(rowValues)) {\n  const r = Number(rowIndex);\n  const values = rowValues[rowIndex];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(r, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date heading and the 25 multiplication problems laid out in\n# the single 20-row x 5-column table (only rows 1, 5, 10, 15, 20 hold the\n# actual \"NNN\u00d7N=\" problems; the other rows are blank spacer rows).\n#\n# Cell addressing (row, col) is used instead of Find/Replace so that the\n# several repeated/aliased values (e.g. one problem's new value equals\n# another problem's old value) can never be mismatched.\n\n$d = $word.ActiveDocument\n\n# 1) Date heading, first paragraph of the document.\n$headingRange = $d.Paragraphs.Item(1).Range\nif ($headingRange.Text.Trim() -eq \"2025-03-02 Sunday\") {\n    $headingRange.Text = \"2025-03-03 Monday\"\n}\n\n# 2) Table of problems.\n$tbl = $d.Tables.Item(1)\n\n$tbl.Cell(1, 1).Range.Text = \"916\u00d79=\"\n$tbl.Cell(1, 2).Range.Text = \"568\u00d74=\"\n$tbl.Cell(1, 3).Range.Text = \"489\u00d77=\"\n$tbl.Cell(1, 4).Range.Text = \"965\u00d78=\"\n$tbl.Cell(1, 5).Range.Text = \"525\u00d78=\"\n\n$tbl.Cell(5, 1).Range.Text = \"366\u00d78=\"\n$tbl.Cell(5, 2).Range.Text = \"964\u00d72=\"\n$tbl.Cell(5, 3).Range.Text = \"544\u00d74=\"\n$tbl.Cell(5, 4).Range.Text = \"985\u00d73=\"\n$tbl.Cell(5, 5).Range.Text = \"724\u00d72=\"\n\n$tbl.Cell(10, 1).Range.Text = \"594\u00d78=\"\n$tbl.Cell(10, 2).Range.Text = \"323\u00d76=\"\n$tbl.Cell(10, 3).Range.Text = \"728\u00d72=\"\n$tbl.Cell(10, 4).Range.Text = \"714\u00d74=\"\n$tbl.Cell(10, 5).Range.Text = \"343\u00d78=\"\n\n$tbl.Cell(15, 1).Range.Text = \"809\u00d77=\"\n$tbl.Cell(15, 2).Range.Text = \"618\u00d74=\"\n$tbl.Cell(15, 3).Range.Text = \"805\u00d75=\"\n$tbl.Cell(15, 4).Range.Text = \"181\u00d74=\"\n$tbl.Cell(15, 5).Range.Text = \"983\u00d79=\"\n\n$tbl.Cell(20, 1).Range.Text = \"800\u00d73=\"\n$tbl.Cell(20, 2).Range.Text = \"866\u00d79=\"\n$tbl.Cell(20, 3).Range.Text = \"746\u00d79=\"\n$tbl.Cell(20, 4).Range.Text = \"731\u00d78=\"\n$tbl.Cell(20, 5).Range.Text = \"201\u00d79=\"\n"}
